$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.439.00'
$ws.Range("E2").Value = '  +1.33%  '

$ws.Range("D3").Value = '2.230.83'
$ws.Range("E3").Value = '  -0.63%  '

$ws.Range("E4").Value = '  +0.53%  '

$ws.Range("E5").Value = '  -1.01%  '

$ws.Range("D6").Value = '''93.93'
$ws.Range("E6").Value = '  -1.80%  '

$ws.Range("E7").Value = '  -1.11%  '

$ws.Range("E8").Value = '  +0.21%  '

$ws.Range("D9").Value = '''0.513'
$ws.Range("E9").Value = '  -2.76%  '

$ws.Range("D10").Value = '''34.52'
$ws.Range("E10").Value = '  -2.07%  '

$ws.Range("D11").Value = '''0.0795'
$ws.Range("E11").Value = '  -2.55%  '

$ws.Range("D12").Value = '''7.10'
$ws.Range("E12").Value = '  -1.91%  '

$ws.Range("E13").Value = '  -0.16%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '2.339.53'
$ws.Range("E14").Value = '  +0.21%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.567.88'
$ws.Range("E15").Value = '  -0.76%  '

$ws.Range("D16").Value = '''0.826'
$ws.Range("E16").Value = '  -1.02%  '

$ws.Range("E17").Value = '  -1.68%  '

$ws.Range("D18").Value = '44.165.83'
$ws.Range("E18").Value = '  +0.23%  '

$ws.Range("D19").Value = '0.0₃0934'
$ws.Range("E19").Value = '  -3.83%  '

$ws.Range("E20").Value = '  -3.83%  '

$ws.Range("D21").Value = '''11.62'
$ws.Range("E21").Value = '  -4.35%  '

$ws.Range("D22").Value = '''64.72'

$ws.Range("D23").Value = '''236.36'
$ws.Range("E23").Value = '  -0.45%  '

$ws.Range("E24").Value = '  -2.28%  '

$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  +0.04%  '

$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").Value = '''1.94'
$ws.Range("E26").Value = '  -3.00%  '

$ws.Range("D27").Value = '''2.31'
$ws.Range("E27").Value = '  +4.63%  '

$ws.Range("D28").Value = '''9.65'
$ws.Range("E28").Value = '  -3.23%  '

$ws.Range("D29").Value = '''37.02'
$ws.Range("E29").Value = '  -1.83%  '

$ws.Range("D30").Value = '''19.73'
$ws.Range("E30").Value = '  -2.34%  '

$ws.Range("D31").Value = '''5.79'
$ws.Range("E31").Value = '  -2.64%  '

$ws.Range("D32").Value = '''148.89'
$ws.Range("E32").Value = '  -2.65%  '

$ws.Range("D33").Value = '''2.61'
$ws.Range("E33").Value = '  +0.22%  '

$ws.Range("E34").Value = '  -2.98%  '

$ws.Range("D35").Value = '''3.09'
$ws.Range("E35").Value = '  -4.38%  '

$ws.Range("E36").Value = '  -0.57%  '

$ws.Range("E37").Value = '  -2.34%  '

$ws.Range("D38").Value = '''1.83'
$ws.Range("E38").Value = '  +4.03%  '

$ws.Range("D39").Value = '''14.68'
$ws.Range("E39").Value = '  +0.62%  '

$ws.Range("D40").Value = '''3.32'
$ws.Range("E40").Value = '  -4.92%  '

$ws.Range("D41").Value = '''3.74'
$ws.Range("E41").Value = '  -3.12%  '

$ws.Range("D42").Value = '''0.0296'
$ws.Range("E42").Value = '  -0.80%  '

$ws.Range("E43").Value = '  +0.29%  '

$ws.Range("D44").Value = '1.818.25'
$ws.Range("E44").Value = '  +3.98%  '

$ws.Range("E45").Value = '  +8.18%  '

$ws.Range("D46").Value = '''78.54'
$ws.Range("E46").Value = '  -5.50%  '

$ws.Range("D47").Value = '''0.186'
$ws.Range("E47").Value = '  -3.31%  '

$ws.Range("D48").Value = '''97.45'
$ws.Range("E48").Value = '  -2.89%  '

$ws.Range("D49").Value = '''4.81'
$ws.Range("E49").Value = '  -2.67%  '

$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D50").Value = '''67.77'
$ws.Range("E50").Value = '  -0.31%  '

$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").Value = '''7.91'
$ws.Range("E51").Value = '  -3.04%  '
